# Apply the vpc.xlsx "networkinterfaces" *interface_type column change.
#
# Summary of the edit:
#   - networkinterfaces sheet/table: a new required column "*interface_type"
#     is inserted between "*instance" and "*subnet"; the two existing NIC
#     rows per VSI are tagged "primary"/"secondary", and the primary NIC
#     names are updated from "vsiNprimarynic" to the new "vsiNnic0"
#     convention.
#   - instances sheet: the *primary_network_interface values are renamed to
#     match ("vsiNprimarynic" -> "vsiNnic0").
#   - the active tab moves back to the first sheet ("vpcheaders").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. networkinterfaces sheet: insert a new "*interface_type" column after
#    "*instance" (so *name, *instance, *interface_type, *subnet,
#    security_groups, floating_ip).
# ---------------------------------------------------------------------------
$wsNic = $wb.Worksheets.Item("networkinterfaces")
$lo = $wsNic.ListObjects.Item(1)

# Grow the table by one column on the right (adds a blank Column6 at F);
# ListColumns.Add(position) appends rather than inserting in this host, so
# we resize first and then move the data/formatting into place ourselves.
$lo.Resize($wsNic.Range("A1:F6"))

# Shift the formatting of the existing *subnet / security_groups /
# floating_ip columns one slot to the right (rightmost first, so a source
# column is never overwritten before it has been copied from).
$wsNic.Range("E1:E6").Copy()
$wsNic.Range("F1:F6").PasteSpecial(-4122) | Out-Null
$wsNic.Range("D1:D6").Copy()
$wsNic.Range("E1:E6").PasteSpecial(-4122) | Out-Null
$wsNic.Range("C1:C6").Copy()
$wsNic.Range("D1:D6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# The new *interface_type column header should look like the other
# required-field headers (e.g. *instance in column B).
$wsNic.Range("B1").Copy()
$wsNic.Range("C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new *interface_type column values first (matches the order in
# which a user would type down the freshly inserted column).
$wsNic.Cells.Item(1, 3).Value2 = "*interface_type"
$wsNic.Cells.Item(2, 3).Value2 = "primary"
$wsNic.Cells.Item(3, 3).Value2 = "secondary"
$wsNic.Cells.Item(4, 3).Value2 = $null
$wsNic.Cells.Item(5, 3).Value2 = "primary"
$wsNic.Cells.Item(6, 3).Value2 = "secondary"

# Re-assert the shifted *subnet / security_groups / floating_ip columns'
# text (Value2; plain Range.Value reads/writes come back as unusable COM
# stubs in this host, so every cell is written individually).
$wsNic.Cells.Item(1, 4).Value2 = "*subnet"
$wsNic.Cells.Item(2, 4).Value2 = "subnet1"
$wsNic.Cells.Item(3, 4).Value2 = "subnet1"
$wsNic.Cells.Item(4, 4).Value2 = $null
$wsNic.Cells.Item(5, 4).Value2 = "subnet2"
$wsNic.Cells.Item(6, 4).Value2 = "subnet2"

$wsNic.Cells.Item(1, 5).Value2 = "security_groups"
$wsNic.Cells.Item(2, 5).Value2 = "sg1"
$wsNic.Cells.Item(3, 5).Value2 = "sg1,sg2"
$wsNic.Cells.Item(4, 5).Value2 = $null
$wsNic.Cells.Item(5, 5).Value2 = "sg1"
$wsNic.Cells.Item(6, 5).Value2 = "sg1,sg2"

$wsNic.Cells.Item(1, 6).Value2 = "floating_ip"
$wsNic.Cells.Item(2, 6).Value2 = "fip1"
$wsNic.Cells.Item(3, 6).Value2 = $null
$wsNic.Cells.Item(4, 6).Value2 = $null
$wsNic.Cells.Item(5, 6).Value2 = $null
$wsNic.Cells.Item(6, 6).Value2 = $null

# Rename the primary-NIC rows to match the new "nic0" naming convention.
$wsNic.Cells.Item(2, 1).Value2 = "vsi1nic0"
$wsNic.Cells.Item(5, 1).Value2 = "vsi2nic0"

# ---------------------------------------------------------------------------
# 2. instances sheet: rename the primary-NIC references to match the new
#    "nic0" naming convention (vsi1primarynic -> vsi1nic0, etc.)
# ---------------------------------------------------------------------------
$wsInstances = $wb.Worksheets.Item("instances")
$wsInstances.Cells.Item(2, 8).Value2 = "vsi1nic0"
$wsInstances.Cells.Item(3, 8).Value2 = "vsi2nic0"

# ---------------------------------------------------------------------------
# 3. Switch the active sheet back to "vpcheaders" (first tab) to match the
#    saved workbook view state.
# ---------------------------------------------------------------------------
$wsHeaders = $wb.Worksheets.Item("vpcheaders")
$wsHeaders.Activate()
Write-Host "Done: interface_type column inserted, nic names updated, active sheet set to vpcheaders."
